# Apply crypto price/volume refresh as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.547.95"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.806.89"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'224.75"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'43.16"
$ws.Range("E8").Value = "  +19.32%  "
$ws.Range("D9").Value = "'0.293"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.0668"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'0.0997"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").Value = "2.065.38"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.810.89"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.630"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.524.49"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'4.41"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'67.33"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'240.65"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'4.37"
$ws.Range("E23").Value = "  +6.85%  "
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").Value = "'170.46"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "'7.68"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'17.44"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'3.80"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'3.88"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'0.0513"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D35").Value = "'87.93"
$ws.Range("E35").Value = "  +8.47%  "
$ws.Range("D36").Value = "'0.650"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "1.318.17"
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'15.01"
$ws.Range("E39").Value = "  +13.83%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "'0.938"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'0.0519"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("D47").Value = "1.966.34"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "'5.81"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'100.60"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("E51").Value = "  +0.70%  "
